$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell ref, new text value.
# Values are written via .Formula with a leading apostrophe so Excel
# always stores them as text (matching the original inlineStr cells)
# instead of inferring a number/date from strings like "211.76".
$updates = @(
    @{ Cell = "D2"; Value = '26.192.11' }
    @{ Cell = "D3"; Value = '1.585.81' }
    @{ Cell = "E3"; Value = '  -0.29%  ' }
    @{ Cell = "E4"; Value = '  -0.16%  ' }
    @{ Cell = "D5"; Value = '211.76' }
    @{ Cell = "E5"; Value = '  +0.86%  ' }
    @{ Cell = "E6"; Value = '  -0.15%  ' }
    @{ Cell = "E7"; Value = '  -0.13%  ' }
    @{ Cell = "E8"; Value = '  -0.39%  ' }
    @{ Cell = "E9"; Value = '  -1.22%  ' }
    @{ Cell = "E10"; Value = '  -1.81%  ' }
    @{ Cell = "D11"; Value = '0.0846' }
    @{ Cell = "E11"; Value = '  +0.16%  ' }
    @{ Cell = "D12"; Value = '1.808.39' }
    @{ Cell = "E12"; Value = '  -0.31%  ' }
    @{ Cell = "D13"; Value = '1.604.76' }
    @{ Cell = "E13"; Value = '  +1.01%  ' }
    @{ Cell = "E14"; Value = '  -1.74%  ' }
    @{ Cell = "E15"; Value = '  -0.28%  ' }
    @{ Cell = "D16"; Value = '63.87' }
    @{ Cell = "E16"; Value = '  -0.99%  ' }
    @{ Cell = "D17"; Value = '26.178.10' }
    @{ Cell = "E17"; Value = '  -0.61%  ' }
    @{ Cell = "E18"; Value = '  -0.58%  ' }
    @{ Cell = "D19"; Value = '213.78' }
    @{ Cell = "E19"; Value = '  +1.36%  ' }
    @{ Cell = "D20"; Value = '7.29' }
    @{ Cell = "E20"; Value = '  -2.03%  ' }
    @{ Cell = "E21"; Value = '  -0.13%  ' }
    @{ Cell = "E22"; Value = '  -0.54%  ' }
    @{ Cell = "D23"; Value = '8.97' }
    @{ Cell = "E23"; Value = '  +0.55%  ' }
    @{ Cell = "E24"; Value = '  -1.79%  ' }
    @{ Cell = "D25"; Value = '144.21' }
    @{ Cell = "E25"; Value = '  -0.46%  ' }
    @{ Cell = "E26"; Value = '  -0.15%  ' }
    @{ Cell = "E27"; Value = '  -1.06%  ' }
    @{ Cell = "E28"; Value = '  -0.85%  ' }
    @{ Cell = "E29"; Value = '  -1.40%  ' }
    @{ Cell = "E30"; Value = '  -2.00%  ' }
    @{ Cell = "D32"; Value = '3.18' }
    @{ Cell = "E32"; Value = '  -0.98%  ' }
    @{ Cell = "D33"; Value = '1.417.75' }
    @{ Cell = "E33"; Value = '  +8.54%  ' }
    @{ Cell = "E34"; Value = '  -1.73%  ' }
    @{ Cell = "E35"; Value = '  -0.48%  ' }
    @{ Cell = "B36"; Value = 'ImmutableX' }
    @{ Cell = "C36"; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = "D36"; Value = '0.586' }
    @{ Cell = "E36"; Value = '  -4.71%  ' }
    @{ Cell = "B37"; Value = 'LidoDAOToken' }
    @{ Cell = "C37"; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = "D37"; Value = '1.45' }
    @{ Cell = "E37"; Value = '  -1.51%  ' }
    @{ Cell = "E38"; Value = '  -1.30%  ' }
    @{ Cell = "D39"; Value = '5.91' }
    @{ Cell = "E39"; Value = '  +5.17%  ' }
    @{ Cell = "E40"; Value = '  +0.94%  ' }
    @{ Cell = "E41"; Value = '  -0.16%  ' }
    @{ Cell = "D42"; Value = '0.949' }
    @{ Cell = "E42"; Value = '  -13.71%  ' }
    @{ Cell = "D43"; Value = '0.765' }
    @{ Cell = "E43"; Value = '  +0.06%  ' }
    @{ Cell = "E44"; Value = '  -0.25%  ' }
    @{ Cell = "D45"; Value = '1.720.27' }
    @{ Cell = "E45"; Value = '  -0.34%  ' }
    @{ Cell = "D46"; Value = '60.87' }
    @{ Cell = "E46"; Value = '  -2.88%  ' }
    @{ Cell = "D47"; Value = '85.25' }
    @{ Cell = "E47"; Value = '  -3.03%  ' }
    @{ Cell = "E48"; Value = '  -0.95%  ' }
    @{ Cell = "E49"; Value = '  -0.89%  ' }
    @{ Cell = "E50"; Value = '  -1.56%  ' }
    @{ Cell = "D51"; Value = '0.999' }
    @{ Cell = "E51"; Value = '  -0.07%  ' }
)

foreach ($u in $updates) {
    $row = [regex]::Match($u.Cell, "\d+").Value
    $target = $ws.Range($u.Cell)
    $target.Formula = "'" + $u.Value
    # Preserve the (unstyled) look of the original cell: column B in this
    # sheet never carries an explicit style for data rows, so copying from
    # it keeps these cells from picking up a stray style index.
    $target.Style = $ws.Range("B" + $row).Style
}
